$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(17, 8).Value = 12045264
$ws.Cells.Item(17, 10).Value = 12045264
$ws.Cells.Item(17, 12).Value = 36135792
$ws.Cells.Item(17, 14).Value = -36136128
$ws.Cells.Item(69, 8).Value = 18500
$ws.Cells.Item(69, 10).Value = 18500
$ws.Cells.Item(69, 12).Value = 55500
$ws.Cells.Item(69, 14).Value = -57248
$ws.Cells.Item(72, 8).Value = 18500
$ws.Cells.Item(72, 10).Value = 18500
$ws.Cells.Item(72, 12).Value = 166500
$ws.Cells.Item(72, 14).Value = -175236
$ws.Cells.Item(80, 8).Value = 1018.625
$ws.Cells.Item(80, 9).Value = 430
$ws.Cells.Item(80, 11).Value = 1290
$ws.Cells.Item(80, 13).Value = -292
$ws.Cells.Item(83, 8).Value = 1018.625
$ws.Cells.Item(83, 9).Value = 430
$ws.Cells.Item(83, 11).Value = 3870
$ws.Cells.Item(83, 13).Value = 1122
$ws.Cells.Item(107, 8).Value = 1255.5769
$ws.Cells.Item(107, 9).Value = 1137.9048
$ws.Cells.Item(107, 10).Value = 1749.8
$ws.Cells.Item(107, 11).Value = 1137.9048
$ws.Cells.Item(107, 12).Value = 1749.8
$ws.Cells.Item(107, 13).Value = 782.0952
$ws.Cells.Item(107, 14).Value = -5589.8
$ws.Cells.Item(112, 8).Value = 2613.7778
$ws.Cells.Item(112, 10).Value = 2708.7058
$ws.Cells.Item(112, 12).Value = 8126.117400000001
$ws.Cells.Item(112, 14).Value = -10342.1174
$ws.Cells.Item(113, 8).Value = 3404.8125
$ws.Cells.Item(113, 9).Value = 2819.7778
$ws.Cells.Item(113, 11).Value = 2819.7778
$ws.Cells.Item(113, 13).Value = 434.2222000000002
$ws.Cells.Item(125, 8).Value = 250.5
$ws.Cells.Item(125, 9).Value = 250.5
$ws.Cells.Item(125, 10).Value = 0
$ws.Cells.Item(125, 11).Value = 2254.5
$ws.Cells.Item(125, 12).Value = 0
$ws.Cells.Item(125, 13).Value = 205.5
$ws.Cells.Item(125, 14).ClearContents()
$ws.Cells.Item(138, 8).Value = 4899.075
$ws.Cells.Item(138, 9).Value = 7538.8667
$ws.Cells.Item(138, 10).Value = 3315.2
$ws.Cells.Item(138, 11).Value = 22616.6001
$ws.Cells.Item(138, 12).Value = 9945.599999999999
$ws.Cells.Item(138, 13).Value = -17476.6001
$ws.Cells.Item(138, 14).Value = -20225.6
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(63, 8).Value = 19694.645
$ws.Cells.Item(63, 9).Value = 3153.2856
$ws.Cells.Item(63, 10).Value = 24519.209
$ws.Cells.Item(63, 11).Value = 3153.2856
$ws.Cells.Item(63, 12).Value = 24519.209
$ws.Cells.Item(63, 13).Value = -2467.2856
$ws.Cells.Item(63, 14).Value = -25891.209
$ws.Cells.Item(66, 8).Value = 19694.645
$ws.Cells.Item(66, 9).Value = 3153.2856
$ws.Cells.Item(66, 10).Value = 24519.209
$ws.Cells.Item(66, 11).Value = 15766.428
$ws.Cells.Item(66, 12).Value = 122596.045
$ws.Cells.Item(66, 13).Value = -12334.428
$ws.Cells.Item(66, 14).Value = -129460.045
$ws.Cells.Item(102, 8).Value = 6417.6
$ws.Cells.Item(102, 9).Value = 6983.4443
$ws.Cells.Item(102, 10).Value = 1325
$ws.Cells.Item(102, 11).Value = 6983.4443
$ws.Cells.Item(102, 12).Value = 1325
$ws.Cells.Item(102, 13).Value = -5361.4443
$ws.Cells.Item(102, 14).Value = -4569
$ws.Cells.Item(110, 8).Value = 1839.6
$ws.Cells.Item(110, 10).Value = 11999
$ws.Cells.Item(110, 12).Value = 11999
$ws.Cells.Item(110, 14).Value = -16089
$ws.Cells.Item(138, 8).Value = 98996.664
$ws.Cells.Item(138, 9).Value = 0
$ws.Cells.Item(138, 10).Value = 98996.664
$ws.Cells.Item(138, 11).Value = 0
$ws.Cells.Item(138, 12).Value = 98996.664
$ws.Cells.Item(138, 13).ClearContents()
$ws.Cells.Item(138, 14).Value = -109276.664
$ws.Cells.Item(139, 8).Value = 76738
$ws.Cells.Item(139, 10).Value = 76738
$ws.Cells.Item(139, 12).Value = 76738
$ws.Cells.Item(139, 14).Value = -87018
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(82, 8).Value = 43792.96
$ws.Cells.Item(82, 9).Value = 10363.2
$ws.Cells.Item(82, 10).Value = 64686.562
$ws.Cells.Item(82, 11).Value = 10363.2
$ws.Cells.Item(82, 12).Value = 64686.562
$ws.Cells.Item(82, 13).Value = -9980.200000000001
$ws.Cells.Item(82, 14).Value = -65452.562
$ws.Cells.Item(85, 8).Value = 43792.96
$ws.Cells.Item(85, 9).Value = 10363.2
$ws.Cells.Item(85, 10).Value = 64686.562
$ws.Cells.Item(85, 11).Value = 10363.2
$ws.Cells.Item(85, 12).Value = 64686.562
$ws.Cells.Item(85, 13).Value = -9037.200000000001
$ws.Cells.Item(85, 14).Value = -67338.56200000001
$ws.Cells.Item(105, 8).Value = 7177.3335
$ws.Cells.Item(105, 9).Value = 8412.294
$ws.Cells.Item(105, 11).Value = 8412.294
$ws.Cells.Item(105, 13).Value = -6665.294
$ws.Cells.Item(107, 8).Value = 13320.211
$ws.Cells.Item(107, 9).Value = 17052.5
$ws.Cells.Item(107, 11).Value = 17052.5
$ws.Cells.Item(107, 13).Value = -15132.5
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 7937.8335
$ws.Cells.Item(16, 9).Value = 5332
$ws.Cells.Item(16, 10).Value = 9799.143
$ws.Cells.Item(16, 11).Value = 5332
$ws.Cells.Item(16, 12).Value = 9799.143
$ws.Cells.Item(16, 13).Value = -5045
$ws.Cells.Item(16, 14).Value = -10373.143
$ws.Cells.Item(68, 8).Value = 70133
$ws.Cells.Item(68, 10).Value = 0
$ws.Cells.Item(68, 12).Value = 0
$ws.Cells.Item(68, 14).ClearContents()
$ws.Cells.Item(71, 8).Value = 70133
$ws.Cells.Item(71, 10).Value = 0
$ws.Cells.Item(71, 12).Value = 0
$ws.Cells.Item(71, 14).ClearContents()
$ws.Cells.Item(74, 8).Value = 45314
$ws.Cells.Item(74, 10).Value = 45314
$ws.Cells.Item(74, 12).Value = 45314
$ws.Cells.Item(74, 14).Value = -47062
$ws.Cells.Item(77, 8).Value = 45314
$ws.Cells.Item(77, 10).Value = 45314
$ws.Cells.Item(77, 12).Value = 135942
$ws.Cells.Item(77, 14).Value = -144678
$ws.Cells.Item(105, 8).Value = 1346.2632
$ws.Cells.Item(105, 9).Value = 1455.6
$ws.Cells.Item(105, 11).Value = 1455.6
$ws.Cells.Item(105, 13).Value = 291.4000000000001
$ws.Cells.Item(107, 8).Value = 1100.4318
$ws.Cells.Item(107, 9).Value = 1061.2812
$ws.Cells.Item(107, 11).Value = 1061.2812
$ws.Cells.Item(107, 13).Value = 858.7188000000001
$ws.Cells.Item(113, 8).Value = 7937.8335
$ws.Cells.Item(113, 9).Value = 5332
$ws.Cells.Item(113, 10).Value = 9799.143
$ws.Cells.Item(113, 11).Value = 5332
$ws.Cells.Item(113, 12).Value = 9799.143
$ws.Cells.Item(113, 13).Value = -3162
$ws.Cells.Item(113, 14).Value = -14139.143
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(4, 8).Value = 23941252
$ws.Cells.Item(4, 9).Value = 43971180
$ws.Cells.Item(4, 10).Value = 96101.48
$ws.Cells.Item(4, 11).Value = 131913540
$ws.Cells.Item(4, 12).Value = 288304.44
$ws.Cells.Item(4, 13).Value = -131913428
$ws.Cells.Item(4, 14).Value = -288528.44
$ws.Cells.Item(17, 8).Value = 499.5
$ws.Cells.Item(17, 9).Value = 300
$ws.Cells.Item(17, 10).Value = 545.53845
$ws.Cells.Item(17, 11).Value = 900
$ws.Cells.Item(17, 12).Value = 1636.61535
$ws.Cells.Item(17, 13).Value = -731
$ws.Cells.Item(17, 14).Value = -1974.61535
$ws.Cells.Item(34, 8).Value = 1803.5714
$ws.Cells.Item(34, 10).Value = 2875
$ws.Cells.Item(34, 12).Value = 8625
$ws.Cells.Item(34, 14).Value = -8793
$ws.Cells.Item(39, 8).Value = 3616
$ws.Cells.Item(39, 10).Value = 3762.9092
$ws.Cells.Item(39, 12).Value = 11288.7276
$ws.Cells.Item(39, 14).Value = -11876.7276
$ws.Cells.Item(75, 8).Value = 4226.375
$ws.Cells.Item(75, 9).Value = 1665.5
$ws.Cells.Item(75, 10).Value = 5080
$ws.Cells.Item(75, 11).Value = 4996.5
$ws.Cells.Item(75, 12).Value = 15240
$ws.Cells.Item(75, 13).Value = -3998.5
$ws.Cells.Item(75, 14).Value = -17236
$ws.Cells.Item(78, 8).Value = 4226.375
$ws.Cells.Item(78, 9).Value = 1665.5
$ws.Cells.Item(78, 10).Value = 5080
$ws.Cells.Item(78, 11).Value = 14989.5
$ws.Cells.Item(78, 12).Value = 45720
$ws.Cells.Item(78, 13).Value = -9997.5
$ws.Cells.Item(78, 14).Value = -55704
$ws.Cells.Item(137, 8).Value = 1707.6666
$ws.Cells.Item(137, 9).Value = 838.5714
$ws.Cells.Item(137, 11).Value = 2515.7142
$ws.Cells.Item(137, 13).Value = 2584.2858
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(107, 8).Value = 47266.863
$ws.Cells.Item(107, 9).Value = 112275.22
$ws.Cells.Item(107, 11).Value = 112275.22
$ws.Cells.Item(107, 13).Value = -110355.22
$ws.Cells.Item(132, 8).Value = 3941497
$ws.Cells.Item(132, 9).Value = 4115.893
$ws.Cells.Item(132, 11).Value = 12347.679
$ws.Cells.Item(132, 13).Value = -9817.679
$ws.Cells.Item(141, 8).Value = 33238.355
$ws.Cells.Item(141, 10).Value = 30212.455
$ws.Cells.Item(141, 12).Value = 30212.455
$ws.Cells.Item(141, 14).Value = -40572.455
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(55, 8).Value = 1174.5333
$ws.Cells.Item(55, 9).Value = 1348.5625
$ws.Cells.Item(55, 10).Value = 975.6429000000001
$ws.Cells.Item(55, 11).Value = 1348.5625
$ws.Cells.Item(55, 12).Value = 975.6429000000001
$ws.Cells.Item(55, 13).Value = -1175.5625
$ws.Cells.Item(55, 14).Value = -1321.6429
$ws.Cells.Item(61, 8).Value = 3014
$ws.Cells.Item(61, 9).Value = 1880.6154
$ws.Cells.Item(61, 10).Value = 6697.5
$ws.Cells.Item(61, 11).Value = 1880.6154
$ws.Cells.Item(61, 12).Value = 6697.5
$ws.Cells.Item(61, 13).Value = -1678.6154
$ws.Cells.Item(61, 14).Value = -7101.5
$ws.Cells.Item(100, 8).Value = 4671.1665
$ws.Cells.Item(100, 9).Value = 3734.875
$ws.Cells.Item(100, 11).Value = 3734.875
$ws.Cells.Item(100, 13).Value = -3193.875
$ws.Cells.Item(113, 8).Value = 3014
$ws.Cells.Item(113, 9).Value = 1880.6154
$ws.Cells.Item(113, 10).Value = 6697.5
$ws.Cells.Item(113, 11).Value = 1880.6154
$ws.Cells.Item(113, 12).Value = 6697.5
$ws.Cells.Item(113, 13).Value = 289.3846000000001
$ws.Cells.Item(113, 14).Value = -11037.5
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(81, 8).Value = 41501.36
$ws.Cells.Item(81, 9).Value = 1715.2106
$ws.Cells.Item(81, 10).Value = 167490.83
$ws.Cells.Item(81, 11).Value = 3430.4212
$ws.Cells.Item(81, 12).Value = 334981.66
$ws.Cells.Item(81, 13).Value = -2369.4212
$ws.Cells.Item(81, 14).Value = -337103.66
$ws.Cells.Item(84, 8).Value = 41501.36
$ws.Cells.Item(84, 9).Value = 1715.2106
$ws.Cells.Item(84, 10).Value = 167490.83
$ws.Cells.Item(84, 11).Value = 17152.106
$ws.Cells.Item(84, 12).Value = 1674908.3
$ws.Cells.Item(84, 13).Value = -11848.106
$ws.Cells.Item(84, 14).Value = -1685516.3
$ws.Cells.Item(100, 8).Value = 1290.2069
$ws.Cells.Item(100, 9).Value = 1245.88
$ws.Cells.Item(100, 11).Value = 2491.76
$ws.Cells.Item(100, 13).Value = -1950.76
$ws.Cells.Item(122, 8).Value = 3209.8
$ws.Cells.Item(122, 9).Value = 2570.8
$ws.Cells.Item(122, 11).Value = 7712.400000000001
$ws.Cells.Item(122, 13).Value = -5262.400000000001
$ws.Cells.Item(132, 8).Value = 1922.95
$ws.Cells.Item(132, 9).Value = 1400.0834
$ws.Cells.Item(132, 11).Value = 4200.2502
$ws.Cells.Item(132, 13).Value = -1670.2502
